$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new experiment rows (25 and 26) at the bottom of the log, matching
# the style already used by the other rows in columns A (short name) and
# N (result note).  We copy the formatting from row 24 (the previous last
# row) with PasteSpecial so the same cell style (grey fill on column A /
# purple fill on column N) is reused instead of creating new styles.
# ---------------------------------------------------------------------------

# Row 25: CNN run
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B25").Value = "PPO use episode distance reward + multiply critic lr + train every episode + self play + shuffle position + CNN + big batch + use cross  vs. Random"
$ws.Range("N24").Copy()
$ws.Range("N25").PasteSpecial(-4122)
$ws.Range("N25").Value = "转圈，估计一个都过不了"
$ws.Range("A25").Value = "ED-SP-Spos-CNN-BB-Crs"

# Row 26: MLP run
$ws.Range("A24").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = "ED-SP-Spos-MLP-BB-Crs"
$ws.Range("B26").Value = "PPO use episode distance reward + multiply critic lr + train every episode + self play + shuffle position + MLP (actor 2 layers) + big batch + use cross  vs. Random"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Update the window/view state: zoom out a bit, stop scrolling the view so
# column A is the leftmost visible column again, and move the active
# selection down onto the newly added data.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 108
$ws.Range("G20").Select()
